$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '94.170.35'
$ws.Range('E2').Value = '  +2.00%  '
$ws.Range('D3').Value = '3.108.05'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.35'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '614.58'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('E7').Value = '  +2.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.390'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.834'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +13.19%  '
$ws.Range('D11').Value = '3.106.99'
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('E12').Value = '  -2.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000245'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -4.90%  '
$ws.Range('D14').Value = '93.698.30'
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.78'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.42'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.98%  '
$ws.Range('D17').Value = '3.686.16'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '3.103.68'
$ws.Range('E18').Value = '  -0.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.67'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.85'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.94'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '445.33'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000202'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.04'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.15'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.37%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.66'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.26'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '85.97'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.36%  '
$ws.Range('D29').Value = '3.270.03'
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.245'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +7.07%  '
$ws.Range('E32').Value = '  +6.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.127'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -7.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.24'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.45%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.96'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.83%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.160'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -7.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.98'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.77%  '
$ws.Range('B39').Value = 'PancakeSwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.90'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.19%  '
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.449'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.31%  '
$ws.Range('B41').Value = 'MantraDAO'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.83'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -6.49%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '477.62'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.71%  '
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '24.02'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +8.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.29'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.06%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.25'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -7.36%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '159.72'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.689'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.86'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -4.41%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.39'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.32'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.16%  '
